$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing data row (2..28) down by one row, working from the
# bottom up so we never overwrite a row before it has been copied. This
# mirrors "insert a row above row 2" without Excel's insert-row format
# guessing (which would otherwise create a spurious extra cell style).
for ($r = 28; $r -ge 2; $r--) {
    $ws.Range("A" + $r + ":S" + $r).Copy()
    $ws.Range("A" + ($r + 1) + ":S" + ($r + 1)).PasteSpecial(-4104)
}

# Column A (the running index 0,1,2,...) uses a bordered/centered style;
# make sure that style is consistent across the whole column after the
# shift above.
$ws.Range("A3").Copy()
$ws.Range("A2:A29").PasteSpecial(-4122)

# Populate the newly freed row 2 with the "紫金山实验室" entry.
$ws.Range("B2").Value = '紫金山实验室'
$ws.Range("C2").Value = '江苏省南京市江宁区'
$ws.Range("D2").Value = '未来网络'
$ws.Range("E2").Value = '网络/区块链/k8s'
$ws.Range("F2").Value = '9:00-18:00'
$ws.Range("G2").Value = '1.5h'
$ws.Range("H2").Value = '基本上很少加班'
$ws.Range("I2").Value = '总包*0.7/12，比例缴纳12%'
$ws.Range("J2").Value = '总包*0.2'
$ws.Range("K2").Value = '试用期6月数；工资不打折'
$ws.Range("L2").Value = '工位大小1.5平方格子，提供联想台式电脑。'
$ws.Range("M2").Value = '按照国家法定节假日。'
$ws.Range("N2").Value = '严格打卡，使用楼下打卡机人脸识别打开。'
$ws.Range("O2").Value = '课题1千万别去，套路太多，管理混乱，领导一言堂，而且是没经验的一言堂'
$ws.Range("Q2").Value = '2022-06-23 09:40:05'
$ws.Range("P2").Value = ""
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""

# Column A is a running index independent of the row shift above; renumber
# it sequentially (0..27) for every data row now that one more row exists.
for ($i = 0; $i -le 27; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
